$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-56 (Generation 0-54) -> Fitness 7310
$ws.Range("C2:C56").Value = 7310

# Rows 57-252 (Generation 55-250) -> Fitness 7293
$ws.Range("C57:C252").Value = 7293
